$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 10000
$ws.Range("I12").Value = 10000
$ws.Range("K12").Value = 10000
$ws.Range("M12").Value = -9830
$ws.Range("H17").Value = 2063.84
$ws.Range("J17").Value = 1785.0476
$ws.Range("L17").Value = 5355.142800000001
$ws.Range("N17").Value = -5691.142800000001
$ws.Range("H28").Value = 6711.875
$ws.Range("I28").Value = 6537.8
$ws.Range("J28").Value = 7002
$ws.Range("K28").Value = 6537.8
$ws.Range("L28").Value = 7002
$ws.Range("M28").Value = -6052.8
$ws.Range("N28").Value = -7972
$ws.Range("H33").Value = 366
$ws.Range("I33").Value = 365.125
$ws.Range("K33").Value = 365.125
$ws.Range("M33").Value = -136.125
$ws.Range("H43").Value = 6199.7144
$ws.Range("I43").Value = 7810.8887
$ws.Range("J43").Value = 3299.6
$ws.Range("K43").Value = 7810.8887
$ws.Range("L43").Value = 3299.6
$ws.Range("M43").Value = -7741.8887
$ws.Range("N43").Value = -3437.6
$ws.Range("H70").Value = 4513.467
$ws.Range("I70").Value = 4256.7144
$ws.Range("K70").Value = 12770.1432
$ws.Range("M70").Value = -12500.1432
$ws.Range("H73").Value = 4513.467
$ws.Range("I73").Value = 4256.7144
$ws.Range("K73").Value = 12770.1432
$ws.Range("M73").Value = -11834.1432
$ws.Range("H88").Value = 6259832.5
$ws.Range("I88").Value = 14291618
$ws.Range("J88").Value = 12887.777
$ws.Range("K88").Value = 14291618
$ws.Range("L88").Value = 12887.777
$ws.Range("M88").Value = -14291212
$ws.Range("N88").Value = -13699.777
$ws.Range("H91").Value = 6259832.5
$ws.Range("I91").Value = 14291618
$ws.Range("J91").Value = 12887.777
$ws.Range("K91").Value = 14291618
$ws.Range("L91").Value = 12887.777
$ws.Range("M91").Value = -14290214
$ws.Range("N91").Value = -15695.777
$ws.Range("H98").Value = 5426.4165
$ws.Range("I98").Value = 5426.4165
$ws.Range("K98").Value = 5426.4165
$ws.Range("M98").Value = -3928.4165
$ws.Range("H111").Value = 1646.4565
$ws.Range("I111").Value = 3311.4167
$ws.Range("J111").Value = 1058.8235
$ws.Range("K111").Value = 9934.250100000001
$ws.Range("L111").Value = 3176.4705
$ws.Range("M111").Value = -6867.250100000001
$ws.Range("N111").Value = -9310.470499999999
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -1746
$ws.Range("H122").Value = 5426.4165
$ws.Range("I122").Value = 5426.4165
$ws.Range("K122").Value = 16279.2495
$ws.Range("M122").Value = -13829.2495
$ws.Range("H132").Value = 1392148
$ws.Range("I132").Value = 3574.386
$ws.Range("J132").Value = 12699104
$ws.Range("K132").Value = 10723.158
$ws.Range("L132").Value = 38097312
$ws.Range("M132").Value = -8193.157999999999
$ws.Range("N132").Value = -38102372
$ws.Range("H135").Value = 1731.1765
$ws.Range("I135").Value = 1702.6
$ws.Range("J135").Value = 1772
$ws.Range("K135").Value = 15323.4
$ws.Range("L135").Value = 15948
$ws.Range("M135").Value = -12788.4
$ws.Range("N135").Value = -21018
$ws.Range("H137").Value = 1517613
$ws.Range("I137").Value = 1726260.9
$ws.Range("K137").Value = 5178782.699999999
$ws.Range("M137").Value = -5176232.699999999
$ws.Range("H138").Value = 4020.8157
$ws.Range("I138").Value = 4205.2
$ws.Range("J138").Value = 3900.5652
$ws.Range("K138").Value = 12615.6
$ws.Range("L138").Value = 11701.6956
$ws.Range("M138").Value = -7475.599999999999
$ws.Range("N138").Value = -21981.6956
$ws.Range("H141").Value = 1664.1562
$ws.Range("I141").Value = 1646.7241
$ws.Range("K141").Value = 4940.1723
$ws.Range("M141").Value = 239.8276999999998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1973.7937
$ws.Range("I32").Value = 2047.8474
$ws.Range("K32").Value = 2047.8474
$ws.Range("M32").Value = -1760.8474
$ws.Range("H45").Value = 30301.4
$ws.Range("I45").Value = 34534.848
$ws.Range("J45").Value = 2784
$ws.Range("K45").Value = 34534.848
$ws.Range("L45").Value = 2784
$ws.Range("M45").Value = -34157.848
$ws.Range("N45").Value = -3538
$ws.Range("H46").Value = 13589.4
$ws.Range("J46").Value = 10665.667
$ws.Range("L46").Value = 10665.667
$ws.Range("N46").Value = -11303.667
$ws.Range("H61").Value = 1999.6842
$ws.Range("I61").Value = 1713.1
$ws.Range("J61").Value = 3074.375
$ws.Range("K61").Value = 1713.1
$ws.Range("L61").Value = 3074.375
$ws.Range("M61").Value = -1501.1
$ws.Range("N61").Value = -3498.375
$ws.Range("H74").Value = 148747.27
$ws.Range("I74").Value = 200521.39
$ws.Range("J74").Value = 3779.7
$ws.Range("K74").Value = 200521.39
$ws.Range("L74").Value = 3779.7
$ws.Range("M74").Value = -199647.39
$ws.Range("N74").Value = -5527.7
$ws.Range("H77").Value = 148747.27
$ws.Range("I77").Value = 200521.39
$ws.Range("J77").Value = 3779.7
$ws.Range("K77").Value = 1002606.95
$ws.Range("L77").Value = 18898.5
$ws.Range("M77").Value = -998238.9500000001
$ws.Range("N77").Value = -27634.5
$ws.Range("H88").Value = 3363.5186
$ws.Range("J88").Value = 3791.1177
$ws.Range("L88").Value = 3791.1177
$ws.Range("N88").Value = -4603.1177
$ws.Range("H91").Value = 3363.5186
$ws.Range("J91").Value = 3791.1177
$ws.Range("L91").Value = 3791.1177
$ws.Range("N91").Value = -6599.1177
$ws.Range("H95").Value = 45000
$ws.Range("J95").Value = 45000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -50492
$ws.Range("H102").Value = 3901.8333
$ws.Range("I102").Value = 2803
$ws.Range("J102").Value = 6099.5
$ws.Range("K102").Value = 2803
$ws.Range("L102").Value = 6099.5
$ws.Range("M102").Value = -1181
$ws.Range("N102").Value = -9343.5
$ws.Range("H132").Value = 9119111
$ws.Range("I132").Value = 2080750.2
$ws.Range("J132").Value = 41671532
$ws.Range("K132").Value = 6242250.6
$ws.Range("L132").Value = 125014596
$ws.Range("M132").Value = -6239720.6
$ws.Range("N132").Value = -125019656
$ws.Range("H136").Value = 1999.6842
$ws.Range("I136").Value = 1713.1
$ws.Range("J136").Value = 3074.375
$ws.Range("K136").Value = 5139.299999999999
$ws.Range("L136").Value = 9223.125
$ws.Range("M136").Value = -2589.299999999999
$ws.Range("N136").Value = -14323.125
$ws.Range("H140").Value = 73200
$ws.Range("J140").Value = 73200
$ws.Range("L140").Value = 73200
$ws.Range("N140").Value = -83560

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 23153030
$ws.Range("I20").Value = 34728052
$ws.Range("J20").Value = 2983.0833
$ws.Range("K20").Value = 34728052
$ws.Range("L20").Value = 2983.0833
$ws.Range("M20").Value = -34727805
$ws.Range("N20").Value = -3477.0833
$ws.Range("H75").Value = 24166
$ws.Range("I75").Value = 16249.5
$ws.Range("J75").Value = 39999
$ws.Range("K75").Value = 16249.5
$ws.Range("L75").Value = 39999
$ws.Range("M75").Value = -15313.5
$ws.Range("N75").Value = -41871
$ws.Range("H78").Value = 24166
$ws.Range("I78").Value = 16249.5
$ws.Range("J78").Value = 39999
$ws.Range("K78").Value = 48748.5
$ws.Range("L78").Value = 119997
$ws.Range("M78").Value = -44068.5
$ws.Range("N78").Value = -129357
$ws.Range("H86").Value = 3401.3333
$ws.Range("I86").Value = 2677.5
$ws.Range("K86").Value = 2677.5
$ws.Range("M86").Value = -1554.5
$ws.Range("H89").Value = 3401.3333
$ws.Range("I89").Value = 2677.5
$ws.Range("K89").Value = 13387.5
$ws.Range("M89").Value = -7771.5
$ws.Range("H99").Value = 87984.25
$ws.Range("I99").Value = 203105
$ws.Range("K99").Value = 203105
$ws.Range("M99").Value = -201607
$ws.Range("H105").Value = 7880418.5
$ws.Range("I105").Value = 527532.7
$ws.Range("K105").Value = 527532.7
$ws.Range("M105").Value = -525785.7
$ws.Range("H107").Value = 2850089.8
$ws.Range("I107").Value = 3345347.8
$ws.Range("J107").Value = 2355
$ws.Range("K107").Value = 3345347.8
$ws.Range("L107").Value = 2355
$ws.Range("M107").Value = -3343427.8
$ws.Range("N107").Value = -6195
$ws.Range("H134").Value = 2686.907
$ws.Range("I134").Value = 2471.0286
$ws.Range("J134").Value = 3631.375
$ws.Range("K134").Value = 7413.085800000001
$ws.Range("L134").Value = 10894.125
$ws.Range("M134").Value = -4878.085800000001
$ws.Range("N134").Value = -15964.125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1555.125
$ws.Range("I16").Value = 1510.6666
$ws.Range("K16").Value = 1510.6666
$ws.Range("M16").Value = -1223.6666
$ws.Range("H20").Value = 100000
$ws.Range("J20").Value = 100000
$ws.Range("L20").Value = 100000
$ws.Range("N20").Value = -100472
$ws.Range("H30").Value = 100000
$ws.Range("J30").Value = 100000
$ws.Range("L30").Value = 100000
$ws.Range("N30").Value = -100182
$ws.Range("H31").Value = 1987759.2
$ws.Range("I31").Value = 2143.139
$ws.Range("J31").Value = 4635247.5
$ws.Range("K31").Value = 2143.139
$ws.Range("L31").Value = 4635247.5
$ws.Range("M31").Value = -1848.139
$ws.Range("N31").Value = -4635837.5
$ws.Range("H34").Value = 1987759.2
$ws.Range("I34").Value = 2143.139
$ws.Range("J34").Value = 4635247.5
$ws.Range("K34").Value = 2143.139
$ws.Range("L34").Value = 4635247.5
$ws.Range("M34").Value = -1941.139
$ws.Range("N34").Value = -4635651.5
$ws.Range("H58").Value = 2640.4
$ws.Range("I58").Value = 1600.8667
$ws.Range("J58").Value = 4199.7
$ws.Range("K58").Value = 1600.8667
$ws.Range("L58").Value = 4199.7
$ws.Range("M58").Value = -1397.8667
$ws.Range("N58").Value = -4605.7
$ws.Range("H59").Value = 80000
$ws.Range("I59").Value = 80000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 80000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -78855
$ws.Range("N59").ClearContents()
$ws.Range("H99").Value = 13191.833
$ws.Range("I99").Value = 13191.833
$ws.Range("K99").Value = 13191.833
$ws.Range("M99").Value = -11693.833
$ws.Range("H105").Value = 2342.7896
$ws.Range("I105").Value = 1913.4667
$ws.Range("K105").Value = 1913.4667
$ws.Range("M105").Value = -166.4666999999999
$ws.Range("H107").Value = 1923775.8
$ws.Range("I107").Value = 3125479.5
$ws.Range("J107").Value = 1049.5
$ws.Range("K107").Value = 3125479.5
$ws.Range("L107").Value = 1049.5
$ws.Range("M107").Value = -3123559.5
$ws.Range("N107").Value = -4889.5
$ws.Range("H113").Value = 1555.125
$ws.Range("I113").Value = 1510.6666
$ws.Range("K113").Value = 1510.6666
$ws.Range("M113").Value = 659.3334
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("N120").ClearContents()
$ws.Range("H121").Value = 40324
$ws.Range("J121").Value = 40324
$ws.Range("L121").Value = 40324
$ws.Range("N121").Value = -42944
$ws.Range("H122").Value = 3282.1
$ws.Range("I122").Value = 1915.5555
$ws.Range("K122").Value = 5746.666499999999
$ws.Range("M122").Value = -3296.666499999999
$ws.Range("H126").Value = 13191.833
$ws.Range("I126").Value = 13191.833
$ws.Range("K126").Value = 39575.499
$ws.Range("M126").Value = -37105.499
$ws.Range("H128").Value = 100000
$ws.Range("J128").Value = 100000
$ws.Range("L128").Value = 100000
$ws.Range("N128").Value = -109960
$ws.Range("H132").Value = 12823448
$ws.Range("I132").Value = 2086.7144
$ws.Range("J132").Value = 27781704
$ws.Range("K132").Value = 6260.1432
$ws.Range("L132").Value = 83345112
$ws.Range("M132").Value = -3730.1432
$ws.Range("N132").Value = -83350172
$ws.Range("H134").Value = 2414.513
$ws.Range("I134").Value = 2327.6
$ws.Range("K134").Value = 6982.799999999999
$ws.Range("M134").Value = -4447.799999999999
$ws.Range("H136").Value = 2640.4
$ws.Range("I136").Value = 1600.8667
$ws.Range("J136").Value = 4199.7
$ws.Range("K136").Value = 4802.6001
$ws.Range("L136").Value = 12599.1
$ws.Range("M136").Value = -2252.6001
$ws.Range("N136").Value = -17699.1

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1235
$ws.Range("J18").Value = 496.5
$ws.Range("L18").Value = 1489.5
$ws.Range("N18").Value = -1827.5
$ws.Range("H50").Value = 312
$ws.Range("I50").Value = 297.5
$ws.Range("J50").Value = 321.66666
$ws.Range("K50").Value = 892.5
$ws.Range("L50").Value = 964.9999799999999
$ws.Range("M50").Value = -411.5
$ws.Range("N50").Value = -1926.99998
$ws.Range("H53").Value = 312
$ws.Range("I53").Value = 297.5
$ws.Range("J53").Value = 321.66666
$ws.Range("K53").Value = 892.5
$ws.Range("L53").Value = 964.9999799999999
$ws.Range("M53").Value = -411.5
$ws.Range("N53").Value = -1926.99998
$ws.Range("H93").Value = 9000
$ws.Range("J93").Value = 9000
$ws.Range("L93").Value = 27000
$ws.Range("N93").Value = -30744
$ws.Range("H113").Value = 1740.4445
$ws.Range("I113").Value = 950
$ws.Range("J113").Value = 1839.25
$ws.Range("K113").Value = 2850
$ws.Range("L113").Value = 5517.75
$ws.Range("M113").Value = -680
$ws.Range("N113").Value = -9857.75
$ws.Range("H122").Value = 586.4815
$ws.Range("I122").Value = 350.53333
$ws.Range("K122").Value = 3154.79997
$ws.Range("M122").Value = -704.79997
$ws.Range("H131").Value = 14634.708
$ws.Range("J131").Value = 1966.25
$ws.Range("L131").Value = 5898.75
$ws.Range("N131").Value = -15978.75
$ws.Range("H134").Value = 2893.1428
$ws.Range("I134").Value = 1637.2
$ws.Range("J134").Value = 6033
$ws.Range("K134").Value = 4911.6
$ws.Range("L134").Value = 18099
$ws.Range("M134").Value = 158.3999999999996
$ws.Range("N134").Value = -28239
$ws.Range("H139").Value = 5358.8237
$ws.Range("I139").Value = 10150.083
$ws.Range("J139").Value = 2745.4092
$ws.Range("K139").Value = 30450.249
$ws.Range("L139").Value = 8236.2276
$ws.Range("M139").Value = -25310.249
$ws.Range("N139").Value = -18516.2276
$ws.Range("H140").Value = 7975.6665
$ws.Range("I140").Value = 8570.933999999999
$ws.Range("K140").Value = 25712.802
$ws.Range("M140").Value = -20532.802
$ws.Range("H141").Value = 10678.277
$ws.Range("I141").Value = 4813.933
$ws.Range("K141").Value = 14441.799
$ws.Range("M141").Value = -9261.798999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 13342769
$ws.Range("I11").Value = 17262300
$ws.Range("K11").Value = 17262300
$ws.Range("M11").Value = -17262161
$ws.Range("H70").Value = 29539782
$ws.Range("I70").Value = 45639664
$ws.Range("K70").Value = 45639664
$ws.Range("M70").Value = -45639394
$ws.Range("H73").Value = 29539782
$ws.Range("I73").Value = 45639664
$ws.Range("K73").Value = 45639664
$ws.Range("M73").Value = -45638728
$ws.Range("H80").Value = 55558820
$ws.Range("J80").Value = 4440.7
$ws.Range("L80").Value = 4440.7
$ws.Range("N80").Value = -6436.7
$ws.Range("H83").Value = 55558820
$ws.Range("J83").Value = 4440.7
$ws.Range("L83").Value = 22203.5
$ws.Range("N83").Value = -32187.5
$ws.Range("H97").Value = 2639.2144
$ws.Range("I97").Value = 1900.1428
$ws.Range("J97").Value = 3378.2856
$ws.Range("K97").Value = 1900.1428
$ws.Range("L97").Value = 3378.2856
$ws.Range("M97").Value = -1404.1428
$ws.Range("N97").Value = -4370.2856
$ws.Range("H113").Value = 3466.6
$ws.Range("I113").Value = 3323.8572
$ws.Range("J113").Value = 3799.6667
$ws.Range("K113").Value = 3323.8572
$ws.Range("L113").Value = 3799.6667
$ws.Range("M113").Value = -1153.8572
$ws.Range("N113").Value = -8139.6667
$ws.Range("H122").Value = 5318.136
$ws.Range("I122").Value = 3160.3333
$ws.Range("K122").Value = 9480.999899999999
$ws.Range("M122").Value = -7030.999899999999
$ws.Range("H126").Value = 26145900
$ws.Range("I126").Value = 1601.3077
$ws.Range("J126").Value = 111114870
$ws.Range("K126").Value = 4803.9231
$ws.Range("L126").Value = 333344610
$ws.Range("M126").Value = -2333.9231
$ws.Range("N126").Value = -333349550
$ws.Range("H132").Value = 2576.6
$ws.Range("I132").Value = 2435.5881
$ws.Range("J132").Value = 3375.6667
$ws.Range("K132").Value = 7306.7643
$ws.Range("L132").Value = 10127.0001
$ws.Range("M132").Value = -4776.7643
$ws.Range("N132").Value = -15187.0001
$ws.Range("H140").Value = 69654.69500000001
$ws.Range("J140").Value = 69654.69500000001
$ws.Range("L140").Value = 69654.69500000001
$ws.Range("N140").Value = -80014.69500000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1539.9286
$ws.Range("I16").Value = 819.5714
$ws.Range("J16").Value = 2260.2856
$ws.Range("K16").Value = 819.5714
$ws.Range("L16").Value = 2260.2856
$ws.Range("M16").Value = -649.5714
$ws.Range("N16").Value = -2600.2856
$ws.Range("H20").Value = 8199.9
$ws.Range("J20").Value = 8199.9
$ws.Range("L20").Value = 8199.9
$ws.Range("N20").Value = -8651.9
$ws.Range("H22").Value = 1780.6666
$ws.Range("J22").Value = 2350
$ws.Range("L22").Value = 2350
$ws.Range("N22").Value = -2940
$ws.Range("H27").Value = 1780.6666
$ws.Range("J27").Value = 2350
$ws.Range("L27").Value = 2350
$ws.Range("N27").Value = -2564
$ws.Range("H100").Value = 3199.4
$ws.Range("I100").Value = 2999.5
$ws.Range("K100").Value = 2999.5
$ws.Range("M100").Value = -2458.5
$ws.Range("H132").Value = 4061.6
$ws.Range("I132").Value = 3128.3635
$ws.Range("J132").Value = 6628
$ws.Range("K132").Value = 9385.0905
$ws.Range("L132").Value = 19884
$ws.Range("M132").Value = -6855.0905
$ws.Range("N132").Value = -24944
$ws.Range("H136").Value = 4398.875
$ws.Range("I136").Value = 4978.7
$ws.Range("J136").Value = 3432.5
$ws.Range("K136").Value = 14936.1
$ws.Range("L136").Value = 10297.5
$ws.Range("M136").Value = -12386.1
$ws.Range("N136").Value = -15397.5
$ws.Range("H138").Value = 79938
$ws.Range("J138").Value = 79938
$ws.Range("L138").Value = 79938
$ws.Range("N138").Value = -90218

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 21451.934
$ws.Range("J45").Value = 25684.166
$ws.Range("L45").Value = 25684.166
$ws.Range("N45").Value = -26666.166
$ws.Range("H62").Value = 6937.25
$ws.Range("J62").Value = 6937.25
$ws.Range("L62").Value = 6937.25
$ws.Range("N62").Value = -8185.25
$ws.Range("H65").Value = 6937.25
$ws.Range("J65").Value = 6937.25
$ws.Range("L65").Value = 34686.25
$ws.Range("N65").Value = -40926.25
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H107").Value = 497
$ws.Range("I107").Value = 415.4074
$ws.Range("J107").Value = 680.5833
$ws.Range("K107").Value = 1246.2222
$ws.Range("L107").Value = 2041.7499
$ws.Range("M107").Value = 673.7778000000001
$ws.Range("N107").Value = -5881.7499
$ws.Range("H113").Value = 774.4074000000001
$ws.Range("I113").Value = 649.1
$ws.Range("J113").Value = 1132.4286
$ws.Range("K113").Value = 1947.3
$ws.Range("L113").Value = 3397.2858
$ws.Range("M113").Value = 222.6999999999998
$ws.Range("N113").Value = -7737.2858
$ws.Range("H126").Value = 1467.5385
$ws.Range("I126").Value = 1458.5
$ws.Range("K126").Value = 4375.5
$ws.Range("M126").Value = -1905.5
$ws.Range("H132").Value = 3963.7307
$ws.Range("I132").Value = 3413.762
$ws.Range("J132").Value = 6273.6
$ws.Range("K132").Value = 10241.286
$ws.Range("L132").Value = 18820.8
$ws.Range("M132").Value = -7711.286
$ws.Range("N132").Value = -23880.8
$ws.Range("H136").Value = 7859.475
$ws.Range("I136").Value = 7794.1387
$ws.Range("J136").Value = 8447.5
$ws.Range("K136").Value = 23382.4161
$ws.Range("L136").Value = 25342.5
$ws.Range("M136").Value = -20832.4161
$ws.Range("N136").Value = -30442.5
